# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 462
$ws1.Range("F4").Value = 7719
$ws1.Range("F5").Value = 89
$ws1.Range("F9").Value = 111
$ws1.Range("F13").Value = 436
$ws1.Range("F14").Value = 62
$ws1.Range("F17").Value = 5571
$ws1.Range("F19").Value = 212
$ws1.Range("F20").Value = 980
$ws1.Range("F22").Value = 319

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 462
$ws4.Range("F4").Value = 7719
$ws4.Range("F5").Value = 89
$ws4.Range("F9").Value = 111
$ws4.Range("F13").Value = 436
$ws4.Range("F14").Value = 62
$ws4.Range("F18").Value = 5571
$ws4.Range("F21").Value = 212
$ws4.Range("F22").Value = 980
$ws4.Range("F24").Value = 319
